$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r2 = New-Object "object[,]" 1,20
$r2[0,0] = "ECs"
$r2[0,1] = "Vegfc"
$r2[0,2] = "Kdr"
$r2[0,3] = "ECs"
$r2[0,4] = 3
$r2[0,5] = 1
$r2[0,6] = 4.425865000000001
$r2[0,7] = 13.277595
$r2[0,8] = 0.4619841037548157
$r2[0,9] = 0.4696223785602887
$r2[0,10] = 3
$r2[0,11] = 1
$r2[0,12] = 153.3847273333333
$r2[0,13] = 460.154182
$r2[0,14] = 0.9432535557163702
$r2[0,15] = 0.9461442014370421
$r2[0,16] = 678.8600962391434
$r2[0,17] = 6109.74086615229
$r2[0,18] = 0.4357681485511704
$r2[0,19] = 0.4443304903398886
$ws.Range("A2:T2").Value = $r2

$r3 = New-Object "object[,]" 1,20
$r3[0,0] = "ECs"
$r3[0,1] = "Vegfc"
$r3[0,2] = "Kdr"
$r3[0,3] = "FAPs"
$r3[0,4] = 3
$r3[0,5] = 1
$r3[0,6] = 4.425865000000001
$r3[0,7] = 13.277595
$r3[0,8] = 0.4619841037548157
$r3[0,9] = 0.4696223785602887
$r3[0,10] = 2
$r3[0,11] = 0.6666666666666666
$r3[0,12] = 0.06306133333333333
$r3[0,13] = 0.189184
$r3[0,14] = 0.0003878014971178633
$r3[0,15] = 0.0003889899333885992
$r3[0,16] = 0.2791009480533334
$r3[0,17] = 2.51190853248
$r3[0,18] = 0.0001791581270807719
$r3[0,19] = 0.0001826783777539622
$ws.Range("A3:T3").Value = $r3

$r4 = New-Object "object[,]" 1,20
$r4[0,0] = "ECs"
$r4[0,1] = "Vegfc"
$r4[0,2] = "Kdr"
$r4[0,3] = "M1"
$r4[0,4] = 3
$r4[0,5] = 1
$r4[0,6] = 4.425865000000001
$r4[0,7] = 13.277595
$r4[0,8] = 0.4619841037548157
$r4[0,9] = 0.4696223785602887
$r4[0,10] = 3
$r4[0,11] = 1
$r4[0,12] = 5.722979333333332
$r4[0,13] = 17.168938
$r4[0,14] = 0.0351939902968738
$r4[0,15] = 0.03530184396657746
$r4[0,16] = 25.32913392712333
$r4[0,17] = 227.96220534411
$r4[0,18] = 0.01625906406485692
$r4[0,19] = 0.01657853593114828
$ws.Range("A4:T4").Value = $r4

$r5 = New-Object "object[,]" 1,20
$r5[0,0] = "ECs"
$r5[0,1] = "Vegfc"
$r5[0,2] = "Kdr"
$r5[0,3] = "M2"
$r5[0,4] = 3
$r5[0,5] = 1
$r5[0,6] = 4.425865000000001
$r5[0,7] = 13.277595
$r5[0,8] = 0.4619841037548157
$r5[0,9] = 0.4696223785602887
$r5[0,10] = 3
$r5[0,11] = 1
$r5[0,12] = 1.951202
$r5[0,13] = 5.853605999999999
$r5[0,14] = 0.01199909701844822
$r5[0,15] = 0.0120358688262385
$r5[0,16] = 8.63575663973
$r5[0,17] = 77.72180975757
$r5[0,18] = 0.00554339208193488
$r5[0,19] = 0.005652313346217754
$ws.Range("A5:T5").Value = $r5

$r6 = New-Object "object[,]" 1,20
$r6[0,0] = "ECs"
$r6[0,1] = "Vegfc"
$r6[0,2] = "Kdr"
$r6[0,3] = "sCs"
$r6[0,4] = 3
$r6[0,5] = 1
$r6[0,6] = 4.425865000000001
$r6[0,7] = 13.277595
$r6[0,8] = 0.4619841037548157
$r6[0,9] = 0.4696223785602887
$r6[0,10] = 2
$r6[0,11] = 1
$r6[0,12] = 1.490433
$r6[0,13] = 2.980866
$r6[0,14] = 0.009165555471189982
$r6[0,15] = 0.006129095836753322
$r6[0,16] = 6.596455249545002
$r6[0,17] = 39.57873149727001
$r6[0,18] = 0.004234340929772751
$r6[0,19] = 0.002878360565280058
$ws.Range("A6:T6").Value = $r6

$r7 = New-Object "object[,]" 1,20
$r7[0,0] = "FAPs"
$r7[0,1] = "Vegfc"
$r7[0,2] = "Kdr"
$r7[0,3] = "ECs"
$r7[0,4] = 3
$r7[0,5] = 1
$r7[0,6] = 4.686805000000001
$r7[0,7] = 14.060415
$r7[0,8] = 0.4892217470254038
$r7[0,9] = 0.4973103589802793
$r7[0,10] = 3
$r7[0,11] = 1
$r7[0,12] = 153.3847273333333
$r7[0,13] = 460.154182
$r7[0,14] = 0.9432535557163702
$r7[0,15] = 0.9461442014370421
$r7[0,16] = 718.8843069895034
$r7[0,17] = 6469.95876290553
$r7[0,18] = 0.4614601524154867
$r7[0,19] = 0.4705273124637651
$ws.Range("A7:T7").Value = $r7

$r8 = New-Object "object[,]" 1,20
$r8[0,0] = "FAPs"
$r8[0,1] = "Vegfc"
$r8[0,2] = "Kdr"
$r8[0,3] = "FAPs"
$r8[0,4] = 3
$r8[0,5] = 1
$r8[0,6] = 4.686805000000001
$r8[0,7] = 14.060415
$r8[0,8] = 0.4892217470254038
$r8[0,9] = 0.4973103589802793
$r8[0,10] = 2
$r8[0,11] = 0.6666666666666666
$r8[0,12] = 0.06306133333333333
$r8[0,13] = 0.189184
$r8[0,14] = 0.0003878014971178633
$r8[0,15] = 0.0003889899333885992
$r8[0,16] = 0.2955561723733334
$r8[0,17] = 2.66000555136
$r8[0,18] = 0.0001897209259190682
$r8[0,19] = 0.0001934487234131992
$ws.Range("A8:T8").Value = $r8

$r9 = New-Object "object[,]" 1,20
$r9[0,0] = "FAPs"
$r9[0,1] = "Vegfc"
$r9[0,2] = "Kdr"
$r9[0,3] = "M1"
$r9[0,4] = 3
$r9[0,5] = 1
$r9[0,6] = 4.686805000000001
$r9[0,7] = 14.060415
$r9[0,8] = 0.4892217470254038
$r9[0,9] = 0.4973103589802793
$r9[0,10] = 3
$r9[0,11] = 1
$r9[0,12] = 5.722979333333332
$r9[0,13] = 17.168938
$r9[0,14] = 0.0351939902968738
$r9[0,15] = 0.03530184396657746
$r9[0,16] = 26.82248815436333
$r9[0,17] = 241.40239338927
$r9[0,18] = 0.01721766541783171
$r9[0,19] = 0.01755597269568444
$ws.Range("A9:T9").Value = $r9

$r10 = New-Object "object[,]" 1,20
$r10[0,0] = "FAPs"
$r10[0,1] = "Vegfc"
$r10[0,2] = "Kdr"
$r10[0,3] = "M2"
$r10[0,4] = 3
$r10[0,5] = 1
$r10[0,6] = 4.686805000000001
$r10[0,7] = 14.060415
$r10[0,8] = 0.4892217470254038
$r10[0,9] = 0.4973103589802793
$r10[0,10] = 3
$r10[0,11] = 1
$r10[0,12] = 1.951202
$r10[0,13] = 5.853605999999999
$r10[0,14] = 0.01199909701844822
$r10[0,15] = 0.0120358688262385
$r10[0,16] = 9.14490328961
$r10[0,17] = 82.30412960649
$r10[0,18] = 0.00587021920609255
$r10[0,19] = 0.005985562246616221
$ws.Range("A10:T10").Value = $r10

$r11 = New-Object "object[,]" 1,20
$r11[0,0] = "FAPs"
$r11[0,1] = "Vegfc"
$r11[0,2] = "Kdr"
$r11[0,3] = "sCs"
$r11[0,4] = 3
$r11[0,5] = 1
$r11[0,6] = 4.686805000000001
$r11[0,7] = 14.060415
$r11[0,8] = 0.4892217470254038
$r11[0,9] = 0.4973103589802793
$r11[0,10] = 2
$r11[0,11] = 1
$r11[0,12] = 1.490433
$r11[0,13] = 2.980866
$r11[0,14] = 0.009165555471189982
$r11[0,15] = 0.006129095836753322
$r11[0,16] = 6.985368836565002
$r11[0,17] = 41.91221301939001
$r11[0,18] = 0.004483989060073811
$r11[0,19] = 0.00304806285080033
$ws.Range("A11:T11").Value = $r11

$r12 = New-Object "object[,]" 1,20
$r12[0,0] = "sCs"
$r12[0,1] = "Vegfc"
$r12[0,2] = "Kdr"
$r12[0,3] = "ECs"
$r12[0,4] = 2
$r12[0,5] = 1
$r12[0,6] = 0.467454
$r12[0,7] = 0.9349080000000001
$r12[0,8] = 0.04879414921978045
$r12[0,9] = 0.03306726245943202
$r12[0,10] = 3
$r12[0,11] = 1
$r12[0,12] = 153.3847273333333
$r12[0,13] = 460.154182
$r12[0,14] = 0.9432535557163702
$r12[0,15] = 0.9461442014370421
$r12[0,16] = 71.70030433087601
$r12[0,17] = 430.2018259852561
$r12[0,18] = 0.04602525474971306
$r12[0,19] = 0.03128639863338839
$ws.Range("A12:T12").Value = $r12

$r13 = New-Object "object[,]" 1,20
$r13[0,0] = "sCs"
$r13[0,1] = "Vegfc"
$r13[0,2] = "Kdr"
$r13[0,3] = "FAPs"
$r13[0,4] = 2
$r13[0,5] = 1
$r13[0,6] = 0.467454
$r13[0,7] = 0.9349080000000001
$r13[0,8] = 0.04879414921978045
$r13[0,9] = 0.03306726245943202
$r13[0,10] = 2
$r13[0,11] = 0.6666666666666666
$r13[0,12] = 0.06306133333333333
$r13[0,13] = 0.189184
$r13[0,14] = 0.0003878014971178633
$r13[0,15] = 0.0003889899333885992
$r13[0,16] = 0.029478272512
$r13[0,17] = 0.176869635072
$r13[0,18] = [double]"1.892244411802328e-05"
$r13[0,19] = [double]"1.286283222143779e-05"
$ws.Range("A13:T13").Value = $r13

$r14 = New-Object "object[,]" 1,20
$r14[0,0] = "sCs"
$r14[0,1] = "Vegfc"
$r14[0,2] = "Kdr"
$r14[0,3] = "M1"
$r14[0,4] = 2
$r14[0,5] = 1
$r14[0,6] = 0.467454
$r14[0,7] = 0.9349080000000001
$r14[0,8] = 0.04879414921978045
$r14[0,9] = 0.03306726245943202
$r14[0,10] = 3
$r14[0,11] = 1
$r14[0,12] = 5.722979333333332
$r14[0,13] = 17.168938
$r14[0,14] = 0.0351939902968738
$r14[0,15] = 0.03530184396657746
$r14[0,16] = 2.675229581283999
$r14[0,17] = 16.051377487704
$r14[0,18] = 0.001717260814185165
$r14[0,19] = 0.001167335339744733
$ws.Range("A14:T14").Value = $r14

$r15 = New-Object "object[,]" 1,20
$r15[0,0] = "sCs"
$r15[0,1] = "Vegfc"
$r15[0,2] = "Kdr"
$r15[0,3] = "M2"
$r15[0,4] = 2
$r15[0,5] = 1
$r15[0,6] = 0.467454
$r15[0,7] = 0.9349080000000001
$r15[0,8] = 0.04879414921978045
$r15[0,9] = 0.03306726245943202
$r15[0,10] = 3
$r15[0,11] = 1
$r15[0,12] = 1.951202
$r15[0,13] = 5.853605999999999
$r15[0,14] = 0.01199909701844822
$r15[0,15] = 0.0120358688262385
$r15[0,16] = 0.9120971797079999
$r15[0,17] = 5.472583078247999
$r15[0,18] = 0.0005854857304207849
$r15[0,19] = 0.0003979932334045245
$ws.Range("A15:T15").Value = $r15

$r16 = New-Object "object[,]" 1,20
$r16[0,0] = "sCs"
$r16[0,1] = "Vegfc"
$r16[0,2] = "Kdr"
$r16[0,3] = "sCs"
$r16[0,4] = 2
$r16[0,5] = 1
$r16[0,6] = 0.467454
$r16[0,7] = 0.9349080000000001
$r16[0,8] = 0.04879414921978045
$r16[0,9] = 0.03306726245943202
$r16[0,10] = 2
$r16[0,11] = 1
$r16[0,12] = 1.490433
$r16[0,13] = 2.980866
$r16[0,14] = 0.009165555471189982
$r16[0,15] = 0.006129095836753322
$r16[0,16] = 0.6967088675820001
$r16[0,17] = 2.786835470328
$r16[0,18] = 0.0004472254813434191
$r16[0,19] = 0.0002026724206729342
$ws.Range("A16:T16").Value = $r16
